# Updates cryptos list values (Price / Volume(1h) columns) to match the
# Mon Aug 28 11:59:19 UTC 2023 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.107.83"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  -0.87%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'1.646.24"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  -1.20%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  -0.69%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'217.55"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  -0.74%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'0.5202"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  -2.55%  "
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'  -0.56%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.2617"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  -1.73%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.06282"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  -1.64%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'20.45"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  -1.67%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.07756"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  -1.22%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'4.472"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  -1.95%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'1.667.28"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  -0.05%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'1.871.20"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  -1.23%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'0.5589"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +0.96%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'0.0₅8006"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  -2.23%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'64.80"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  -1.63%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'26.096.66"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  -0.97%  "
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = "'  -0.61%  "
$ws.Range("E19").ClearFormats()
$ws.Range("E20").Value = "'  -0.60%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'192.45"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  -0.57%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'10.11"
$ws.Range("D22").ClearFormats()
$ws.Range("D23").Value = "'5.950"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  -1.36%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'1.005"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  -0.74%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'146.18"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").Value = "'0.1200"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  -2.28%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'7.174"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  -0.37%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D29").Value = "'1.470"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  -2.01%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'0.05616"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  -4.00%  "
$ws.Range("E30").ClearFormats()
$ws.Range("E31").Value = "'  -1.57%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'3.455"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  -3.66%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'3.344"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +2.02%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'1.594"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  -0.57%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'2.789"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  -1.44%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D37").Value = "'0.9364"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  -3.40%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'0.5674"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -2.58%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'5.970"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +2.41%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'0.01583"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").Value = "'2.569"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  -1.36%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'1.049.85"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -1.36%  "
$ws.Range("E42").ClearFormats()
$ws.Range("E43").Value = "'  -0.67%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.8416"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  -2.28%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'102.28"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  -2.27%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'1.782.27"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  -1.29%  "
$ws.Range("E46").ClearFormats()
$ws.Range("E47").Value = "'  +2.91%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'57.25"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  -0.92%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'1.007"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  -0.64%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.05335"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  +3.27%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'0.4334"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  -1.35%  "
$ws.Range("E51").ClearFormats()
